$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the hyperlinks that must survive the edit (B4, B5) before we
# touch anything, since this runtime's Hyperlinks.Delete() on a range
# clears the whole worksheet's hyperlink collection.
$b4 = $ws.Range("B4")
$b5 = $ws.Range("B5")
$b4Text = $b4.Value2
$b5Text = $b5.Value2

# Drop all hyperlinks (the only reliable deletion mechanism here); the
# three supplier-link cells whose URLs changed (B2, B3, B7) should end up
# with plain text and no hyperlink, while B4 and B5 get their hyperlink
# recreated below.
$ws.Range("B2").Hyperlinks.Delete()

# Update the three component source URLs that changed.
$ws.Range("B7").Value = "http://china.rs-online.com/web/p/surface-mount-fixed-resistors/6791263/"
$ws.Range("B3").Value = "http://china.rs-online.com/web/p/surface-mount-fixed-resistors/2230562/"
$ws.Range("B2").Value = "http://china.rs-online.com/web/p/ceramic-multilayer-capacitors/7236054/"

# Recreate the hyperlinks that should remain (B4, B5), and restore their
# original cell formatting.
$ws.Hyperlinks.Add($b4, $b4Text) | Out-Null
$ws.Hyperlinks.Add($b5, $b5Text) | Out-Null

$b4.Style = "Hyperlink"
$b4.WrapText = $true
$b5.Style = "Hyperlink"

# The rows whose wrapped text changed got a shorter recalculated row
# height when Excel re-rendered them.
$ws.Rows.Item(2).RowHeight = 24.75
$ws.Rows.Item(3).RowHeight = 24.75
$ws.Rows.Item(7).RowHeight = 24.75

# Move the active selection to E11, as recorded in the saved view state.
$ws.Range("E11").Select()
